# "fixed name_check_in and number_check_in"
# Adds a new date column (July 18 2016) with a check-in time for row 2,
# marks js (row 4) as PAID, and adds a new attendee row (ccff).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Attendance")

# New column F: check-in date header
$ws.Range("F1").Value = "July 18 2016"

# Row 2 (xc): check-in time for the new date column
$ws.Range("F2").Value = "02:34 PM"

# Row 4 (js): mark as Paid
$ws.Range("C4").Value = "PAID"

# Row 5: new attendee entry
$ws.Range("A5").Value = 0
$ws.Range("B5").Value = "ccff"

$wb.Save()
